$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.301.66'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.081.66'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.17%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '328.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9996'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5220'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.80%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4318'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08828'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '46.72'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.23%  '
$ws.Range('E11').Value = '  +2.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.45'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.084.72'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.23%  '
$ws.Range('E14').Value = '  +1.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.665'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '95.45'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.05%  '
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06633'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.87'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.305'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.354.91'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.300'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.332.35'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.594'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.03'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '131.35'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.191'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1070'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.650'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +20.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.181'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.861'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.933'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02571'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06682'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.69'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.26%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.451'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2263'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6816'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.246'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9997'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('E45').Value = '  +2.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6366'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.202'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.612'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.252'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.185'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '81.56'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.57%  '
